# Updated symbol list on Mon Dec 12 19:47:37 UTC 2022 with GitHub Actions
#
# Applies the per-row "Price" (column D) refresh, a couple of "Volume(1h)"
# (column E) label tweaks, and the CEJI / KickToken / BKEXToken row content
# rotation (columns B-E for rows 41-43).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$addr,
        [string]$value
    )
    # A leading apostrophe forces Excel to store the value as literal text
    # even when it looks numeric (e.g. "273.77"), matching the workbook's
    # existing inline-string cell type instead of converting it to a
    # floating point number. Resetting the style back to Normal afterwards
    # drops the "quote prefix" formatting flag Excel applies so the cell's
    # style stays identical to the original.
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

# ---- Column D "Price" refreshes ----
Set-TextValue "D2"  "273.77"
Set-TextValue "D3"  "21.13"
Set-TextValue "D4"  "6.208"
Set-TextValue "D5"  "0.06175"
Set-TextValue "D6"  "3.575"
Set-TextValue "D8"  "6.529"
Set-TextValue "D9"  "0.8231"
Set-TextValue "D11" "0.08248"
Set-TextValue "D12" "0.03437"
Set-TextValue "D13" "0.03150"
Set-TextValue "D15" "3.763"
Set-TextValue "D16" "0.001621"
Set-TextValue "D17" "0.04697"
Set-TextValue "D19" "0.006135"
Set-TextValue "D21" "0.0001501"
Set-TextValue "D22" "3.721"
Set-TextValue "D25" "0.3339"
Set-TextValue "D26" "0.1232"
Set-TextValue "D28" "0.0002738"
Set-TextValue "D40" "0.04753"
Set-TextValue "D44" "0.01034"
Set-TextValue "D45" "0.00006565"
Set-TextValue "D47" "0.7233"
Set-TextValue "D49" "0.00001901"

# ---- Column E "Volume(1h)" label tweaks ----
$ws.Range("E19").Value = "18HotbitTokenHTBBestin24h"

# ---- Rows 41-43: CEJI / KickToken / BKEXToken rotate into each other ----
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D41" "0.007039"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1106"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.003521"
$ws.Range("E43").Value = "42CEJICEJI"
